$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-31 20:50:57"

# Determine the last used row on the sheet (data starts at row 2, header at row 1)
$lastRow = $ws.UsedRange.Rows.Count

# Update the "timestamp" column (O) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Update the "productAriaLabel" column (M) for the two products that became
# unavailable online, inserting " - Online kein Bestand" right before the price.
$ws.Range("M378").Value = "Betty Bossi Naturaplan Bio Mungbohnen-Sprossen - Online kein Bestand 2.20 Schweizer Franken"
$ws.Range("M518").Value = "Yolo Crunchies - Online kein Bestand 6.95 Schweizer Franken"

$wb.Save()
